$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0374241828918457
$ws.Range("C2").Value = 0.05635991096496582
$ws.Range("D2").Value = 0.01432099342346191
$ws.Range("E2").Value = 0.0374549388885498
$ws.Range("F2").Value = 0.01039361953735352
$ws.Range("G2").Value = 0.1382022857666016
$ws.Range("H2").Value = 0.04408912658691407
$ws.Range("I2").Value = 0.07538542747497559
$ws.Range("J2").Value = 0.02771964073181152
$ws.Range("K2").Value = 0.06352519989013672
$ws.Range("L2").Value = 0.01035785675048828
$ws.Range("M2").Value = 0.07379536628723145
$ws.Range("B3").Value = 0.2319737911224365
$ws.Range("C3").Value = 0.08908252716064453
$ws.Range("D3").Value = 0.03519744873046875
$ws.Range("E3").Value = 0.01848273277282715
$ws.Range("F3").Value = 0.01509771347045898
$ws.Range("G3").Value = 0.01677374839782715
$ws.Range("H3").Value = 0.2128887176513672
$ws.Range("I3").Value = 0.06731266975402832
$ws.Range("J3").Value = 0.1457373142242432
$ws.Range("K3").Value = 0.04217424392700195
$ws.Range("L3").Value = 0.04871301651000977
$ws.Range("M3").Value = 0.02279915809631348
$ws.Range("B4").Value = 0.05716004371643067
$ws.Range("C4").Value = 0.02171125411987305
$ws.Range("D4").Value = 0.03658490180969239
$ws.Range("E4").Value = 0.04089937210083008
$ws.Range("F4").Value = 0.2064003467559815
$ws.Range("G4").Value = 0.02874436378479004
$ws.Range("H4").Value = 0.05865049362182617
$ws.Range("I4").Value = 0.06721343994140624
$ws.Range("J4").Value = 0.04904403686523438
$ws.Range("K4").Value = 0.06456432342529297
$ws.Range("L4").Value = 0.07051692008972169
$ws.Range("M4").Value = 0.02154617309570312
$ws.Range("B5").Value = 0.0354374885559082
$ws.Range("C5").Value = 0.03339834213256836
$ws.Range("D5").Value = 0.03480215072631836
$ws.Range("E5").Value = 0.0316965103149414
$ws.Range("H5").Value = 0.03844742774963379
$ws.Range("I5").Value = 0.03392887115478516
$ws.Range("J5").Value = 0.02048883438110351
$ws.Range("K5").Value = 0.03802752494812012
$ws.Range("B6").Value = 0.9877229690551758
$ws.Range("C6").Value = 0.04915952682495117
$ws.Range("D6").Value = 0.9749974727630615
$ws.Range("E6").Value = 0.04991369247436524
$ws.Range("F6").Value = 2.18835391998291
$ws.Range("G6").Value = 0.04042840003967285
$ws.Range("H6").Value = 0.5025713443756104
$ws.Range("I6").Value = 0.03336911201477051
$ws.Range("J6").Value = 1.012448644638062
$ws.Range("K6").Value = 0.03944311141967773
$ws.Range("L6").Value = 0.7677135467529297
$ws.Range("M6").Value = 0.02432136535644531
$ws.Range("B7").Value = 1.032159948348999
$ws.Range("C7").Value = 0.1392716884613037
$ws.Range("D7").Value = 0.5096531867980957
$ws.Range("E7").Value = 0.06172366142272949
$ws.Range("F7").Value = 0.7097468376159668
$ws.Range("G7").Value = 0.03217315673828125
$ws.Range("H7").Value = 1.149969005584717
$ws.Range("I7").Value = 0.1227193355560303
$ws.Range("J7").Value = 0.5003772735595703
$ws.Range("K7").Value = 0.07103452682495118
$ws.Range("L7").Value = 0.7030947208404541
$ws.Range("M7").Value = 0.03682327270507812
